$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B230").Value = 7482832
$ws.Range("F230").Value = "Barcelona Guayaquil"
$ws.Range("G230").Value = "Guayaquil City"
$ws.Range("H230").Value = 2
$ws.Range("I230").Value = 1
$ws.Range("J230").Value = "H"
$ws.Range("K230").Value = 1.363
$ws.Range("L230").Value = 5
$ws.Range("M230").Value = 7.5
$ws.Range("N230").Value = 1.444
$ws.Range("O230").Value = 4
$ws.Range("P230").Value = 8
$ws.Range("Q230").Value = -1.25
$ws.Range("R230").Value = 2.05
$ws.Range("S230").Value = 1.75
$ws.Range("U230").Value = 1.95
$ws.Range("V230").Value = 1.85
$ws.Range("W230").Value = 0.444
$ws.Range("Y230").Value = -1
$ws.Range("Z230").Value = -0.5
$ws.Range("AA230").Value = 0.375
$ws.Range("AB230").Value = 0.95
$ws.Range("B231").Value = 7483306
$ws.Range("F231").Value = "Tecnico Universitario"
$ws.Range("G231").Value = "Club Atletico Libertad"
$ws.Range("H231").Value = 1
$ws.Range("I231").Value = 1
$ws.Range("J231").Value = "D"
$ws.Range("K231").Value = 1.5
$ws.Range("L231").Value = 4.333
$ws.Range("M231").Value = 5.75
$ws.Range("N231").Value = 1.533
$ws.Range("O231").Value = 4.2
$ws.Range("P231").Value = 5.5
$ws.Range("Q231").Value = -1
$ws.Range("R231").Value = 1.925
$ws.Range("S231").Value = 1.875
$ws.Range("T231").Value = 2.25
$ws.Range("U231").Value = 1.8
$ws.Range("V231").Value = 2
$ws.Range("X231").Value = 3.2
$ws.Range("Y231").Value = -1
$ws.Range("AA231").Value = 0.875
$ws.Range("AB231").Value = -0.5
$ws.Range("AC231").Value = 0.5
$ws.Range("B232").Value = 7483188
$ws.Range("F232").Value = "Gualaceo SC"
$ws.Range("G232").Value = "Emelec"
$ws.Range("H232").Value = 0
$ws.Range("I232").Value = 2
$ws.Range("J232").Value = "A"
$ws.Range("K232").Value = 3.6
$ws.Range("L232").Value = 3.3
$ws.Range("M232").Value = 2.05
$ws.Range("N232").Value = 2.6
$ws.Range("O232").Value = 3.25
$ws.Range("P232").Value = 2.75
$ws.Range("Q232").Value = 0
$ws.Range("R232").Value = 1.8
$ws.Range("S232").Value = 2
$ws.Range("U232").Value = 1.975
$ws.Range("V232").Value = 1.825
$ws.Range("W232").Value = -1
$ws.Range("Y232").Value = 1.75
$ws.Range("Z232").Value = -1
$ws.Range("AA232").Value = 1
$ws.Range("AB232").Value = -1
$ws.Range("AC232").Value = 0.825
$ws.Range("B233").Value = 7482867
$ws.Range("F233").Value = "Cumbaya FC"
$ws.Range("G233").Value = "LDU Quito"
$ws.Range("I233").Value = 2
$ws.Range("J233").Value = "A"
$ws.Range("K233").Value = 5.25
$ws.Range("L233").Value = 3.75
$ws.Range("M233").Value = 1.65
$ws.Range("N233").Value = 9
$ws.Range("O233").Value = 4.5
$ws.Range("P233").Value = 1.363
$ws.Range("Q233").Value = 1.25
$ws.Range("R233").Value = 1.975
$ws.Range("S233").Value = 1.825
$ws.Range("T233").Value = 2.5
$ws.Range("U233").Value = 1.825
$ws.Range("V233").Value = 1.975
$ws.Range("X233").Value = -1
$ws.Range("Y233").Value = 0.363
$ws.Range("Z233").Value = 0.4875
$ws.Range("AA233").Value = -0.5
$ws.Range("AB233").Value = 0.825
$ws.Range("AC233").Value = -1
$ws.Range("B238").Value = 7528857
$ws.Range("F238").Value = "Universidad Catolica del Ecuador"
$ws.Range("G238").Value = "Barcelona Guayaquil"
$ws.Range("H238").Value = 0
$ws.Range("I238").Value = 1
$ws.Range("K238").Value = 1.533
$ws.Range("L238").Value = 4
$ws.Range("M238").Value = 5.5
$ws.Range("N238").Value = 1.5
$ws.Range("O238").Value = 4.333
$ws.Range("P238").Value = 5.25
$ws.Range("Q238").Value = -1
$ws.Range("R238").Value = 1.8
$ws.Range("S238").Value = 2
$ws.Range("T238").Value = 3
$ws.Range("U238").Value = 1.975
$ws.Range("V238").Value = 1.825
$ws.Range("Y238").Value = 4.25
$ws.Range("AA238").Value = 1
$ws.Range("AB238").Value = -1
$ws.Range("AC238").Value = 0.825
$ws.Range("B239").Value = 7528848
$ws.Range("F239").Value = "Emelec"
$ws.Range("G239").Value = "Deportivo Cuenca"
$ws.Range("I239").Value = 1
$ws.Range("J239").Value = "H"
$ws.Range("K239").Value = 1.75
$ws.Range("L239").Value = 3.5
$ws.Range("M239").Value = 4.2
$ws.Range("N239").Value = 2.4
$ws.Range("O239").Value = 3.1
$ws.Range("P239").Value = 2.75
$ws.Range("R239").Value = 2.05
$ws.Range("S239").Value = 1.75
$ws.Range("U239").Value = 1.8
$ws.Range("V239").Value = 2
$ws.Range("W239").Value = 1.4
$ws.Range("X239").Value = -1
$ws.Range("Z239").Value = 1.05
$ws.Range("AA239").Value = -1
$ws.Range("AB239").Value = 0.8
$ws.Range("B240").Value = 7528852
$ws.Range("F240").Value = "Delfin SC"
$ws.Range("G240").Value = "Tecnico Universitario"
$ws.Range("I240").Value = 2
$ws.Range("J240").Value = "D"
$ws.Range("K240").Value = 2.1
$ws.Range("L240").Value = 3.4
$ws.Range("M240").Value = 3.1
$ws.Range("N240").Value = 2.1
$ws.Range("O240").Value = 3.4
$ws.Range("P240").Value = 3.1
$ws.Range("R240").Value = 1.8
$ws.Range("S240").Value = 2
$ws.Range("U240").Value = 1.9
$ws.Range("V240").Value = 1.9
$ws.Range("W240").Value = -1
$ws.Range("X240").Value = 2.4
$ws.Range("Z240").Value = -0.5
$ws.Range("AA240").Value = 0.5
$ws.Range("AB240").Value = 0.8999999999999999
$ws.Range("B241").Value = 7528858
$ws.Range("F241").Value = "Orense"
$ws.Range("G241").Value = "SD Aucas"
$ws.Range("H241").Value = 1
$ws.Range("I241").Value = 2
$ws.Range("K241").Value = 2.2
$ws.Range("L241").Value = 3.2
$ws.Range("M241").Value = 3.2
$ws.Range("N241").Value = 1.95
$ws.Range("O241").Value = 3.2
$ws.Range("P241").Value = 3.8
$ws.Range("Q241").Value = -0.5
$ws.Range("R241").Value = 1.95
$ws.Range("S241").Value = 1.85
$ws.Range("T241").Value = 2.25
$ws.Range("U241").Value = 1.85
$ws.Range("V241").Value = 1.95
$ws.Range("Y241").Value = 2.8
$ws.Range("AA241").Value = 0.8500000000000001
$ws.Range("AB241").Value = 0.8500000000000001
$ws.Range("AC241").Value = -1
